$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store"
$ws.Range("B2").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\bin\.DS_Store"
$ws.Range("A3").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store"
$ws.Range("B3").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\bin\.DS_Store"
$ws.Range("A4").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store"
$ws.Range("B4").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\Docs\.DS_Store"
$ws.Range("A5").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\.DS_Store"
$ws.Range("B5").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\.tar\opt\ddn\others\RAW_BM\XDD\xdd_r4563_20150503\xdd\Docs\.DS_Store"
$ws.Range("A6").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\08456059.pdf"
$ws.Range("B6").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\08456059.pdf"
$ws.Range("A7").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\20180006549.pdf"
$ws.Range("B7").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\20180006549.pdf"
$ws.Range("A8").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\235434231_RAW_PAN.pix"
$ws.Range("B8").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\235434231_RAW_PAN.pix"
$ws.Range("A9").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\235434251_RAW_PAN_LOG.txt"
$ws.Range("B9").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\235434251_RAW_PAN_LOG.txt"
$ws.Range("A10").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\3.205-CVE_to_be_installed.txt"
$ws.Range("B10").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\3.205-CVE_to_be_installed.txt"
$ws.Range("A11").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\3.207-CVE_to_be_installed.txt"
$ws.Range("B11").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\3.207-CVE_to_be_installed.txt"
$ws.Range("A12").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\anaconda-ks.cfg"
$ws.Range("B12").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\anaconda-ks.cfg"
$ws.Range("A13").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\authorized_keys-4-8-17"
$ws.Range("B13").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\authorized_keys-4-8-17"
$ws.Range("A14").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\BAND_RPC.txt"
$ws.Range("B14").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\BAND_RPC.txt"
$ws.Range("A15").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\buttons.print.min.js"
$ws.Range("B15").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\buttons.print.min.js"
$ws.Range("A16").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\catalina.2023-07-14.log"
$ws.Range("B16").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\catalina.2023-07-14.log"
$ws.Range("A17").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\catalina.2023-07-16.log"
$ws.Range("B17").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\catalina.2023-07-16.log"
$ws.Range("A18").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\catalina.2023-07-23.log"
$ws.Range("B18").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\catalina.2023-07-23.log"
$ws.Range("A19").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\catalina.2023-08-26.log"
$ws.Range("B19").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\catalina.2023-08-26.log"
$ws.Range("A20").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\covid_excel_20-38.xlsx"
$ws.Range("B20").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\covid_excel_20-38.xlsx"
$ws.Range("A21").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\ddn-restapi-core-gs-4.2.1-3.3d53e4_dirty.noarch.rpm"
$ws.Range("B21").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\ddn-restapi-core-gs-4.2.1-3.3d53e4_dirty.noarch.rpm"
$ws.Range("A22").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\EaseUS Partition Master 9.2.2 Technician Edition (FULL + Patch).zip"
$ws.Range("B22").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\EaseUS Partition Master 9.2.2 Technician Edition (FULL + Patch).zip"
$ws.Range("A23").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\fio-3.1-1.el6.x86_64.rpm"
$ws.Range("B23").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\fio-3.1-1.el6.x86_64.rpm"
$ws.Range("A24").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\gpfs.gss.pmsensors-4.2.3-5.el7.x86_64.rpm"
$ws.Range("B24").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\gpfs.gss.pmsensors-4.2.3-5.el7.x86_64.rpm"
$ws.Range("A25").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\gwm66c08.jpg"
$ws.Range("B25").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\gwm66c08.jpg"
$ws.Range("A26").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\import CSV.txt"
$ws.Range("B26").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\import CSV.txt"
$ws.Range("A27").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\KKD geo tagging (1).xlsx"
$ws.Range("B27").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\KKD geo tagging (1).xlsx"
$ws.Range("A28").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\LS____.pdf"
$ws.Range("B28").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\LS____.pdf"
$ws.Range("A29").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\mmhealth.sh"
$ws.Range("B29").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\mmhealth.sh"
$ws.Range("A30").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\mounted"
$ws.Range("B30").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\mounted"
$ws.Range("A31").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\NewGCPs.bk"
$ws.Range("B31").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\NewGCPs.bk"
$ws.Range("A32").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\SR159360_ddn_showall_20210222-13_54_51.tar.gz"
$ws.Range("B32").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\SR159360_ddn_showall_20210222-13_54_51.tar.gz"
$ws.Range("A33").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\SR170158_ddn_showall_20210902-13_20_38.tar.gz"
$ws.Range("B33").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\SR170158_ddn_showall_20210902-13_20_38.tar.gz"
$ws.Range("A34").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\syslog"
$ws.Range("B34").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\syslog"
$ws.Range("A35").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-1\table.sql"
$ws.Range("B35").Value = "D:\adhvik\adh\Hackathon\space hack\Data RR\data Set\topic12\dataset1\folder-2\table.sql"
